# Insert a new "Audience" slide as the 5th slide of the deck (right after
# "Data requirements" and before "Data DESCRIPTION"), using the same
# "Title and Content" layout as the other text slides.

$p = $ppt.ActivePresentation

# Custom layout #2 on the slide master is "Title and Content" - the layout
# already used by the surrounding slides (Background, Problem statement,
# Data requirements, Data DESCRIPTION, Scope, references, ...).
$layout = $p.SlideMaster.CustomLayouts.Item(2)

# Add the new slide at position 5.
$s = $p.Slides.AddSlide(5, $layout)

# Title placeholder.
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Audience"

# Body / content placeholder (idx=13 "quarter" placeholder from the layout).
$body = $s.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "The target audience for this project would be the management for the waffle company who are interested in opening their franchise in NA`rThis paper will also interest students of the datascience field as a reference "
